$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A values (location names) - same 7 names/order as before the edit
$colA = @(
  "울부짖는 참나무",
  "세나리온 자치령",
  "장인의 정원",
  "전사의 정원",
  "신전 정원",
  "상인의 정원",
  "달의 신전"
)

# Column B values (descriptions) - newly filled in for every row
$colB = @(
  "피난처",
  "자치령",
  "전문가들에게 기술을 배울수 있음",
  "다르나서스 입구가 있음",
  "은행과 텔드랏실 아래로 가는 포털이 있음",
  "상인들에게 물건을 살 수 있음",
  "신전"
)

for ($i = 0; $i -lt $colA.Length; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $colA[$i]
  $ws.Cells.Item($row, 2).Value = $colB[$i]
}

# Widen column A to fit the longer location names
$ws.Columns.Item(1).ColumnWidth = 25.4

# Move the active selection from B2 to A2
$ws.Range("A2").Select()
